$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are set directly. Cells whose new text would otherwise be
# auto-parsed by Excel as a number are entered with a leading apostrophe
# (standard Excel "store as text" input) so they stay text, matching the
# original inline-string cell type.

$ws.Range("D2").Value = "27.383.41"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.778.41"
$ws.Range("E3").Value = "  +3.72%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'313.75"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'0.5362"
$ws.Range("E7").Value = "  +13.20%  "
$ws.Range("D8").Value = "'0.3773"
$ws.Range("E8").Value = "  +9.22%  "
$ws.Range("D9").Value = "'42.82"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").Value = "'0.07401"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").Value = "'1.093"
$ws.Range("E11").Value = "  +5.40%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'20.67"
$ws.Range("E13").Value = "  +4.63%  "
$ws.Range("D14").Value = "'6.098"
$ws.Range("E14").Value = "  +4.84%  "
$ws.Range("D15").Value = "1.781.85"
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("D16").Value = "'6.992"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "'89.60"
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("D18").Value = "'0.00001055"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").Value = "'0.06432"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("D21").Value = "'16.78"
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("E22").Value = "  +5.42%  "
$ws.Range("D23").Value = "27.422.51"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").Value = "'11.19"
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("D25").Value = "'2.092"
$ws.Range("D26").Value = "'155.44"
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").Value = "'20.20"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").Value = "'2.374"
$ws.Range("E28").Value = "  +15.17%  "
$ws.Range("D29").Value = "1.988.21"
$ws.Range("E29").Value = "  +3.84%  "
$ws.Range("D30").Value = "'121.28"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'1.078"
$ws.Range("E31").Value = "  +5.07%  "
$ws.Range("D32").Value = "'0.1032"
$ws.Range("E32").Value = "  +12.74%  "
$ws.Range("D33").Value = "'5.588"
$ws.Range("E33").Value = "  +5.57%  "
$ws.Range("D34").Value = "'3.622"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "'0.02261"
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("D36").Value = "'0.05965"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "'4.917"
$ws.Range("E37").Value = "  +4.75%  "
$ws.Range("D38").Value = "'11.27"
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("D44").Value = "'13.25"
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("D45").Value = "'0.5779"
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("D46").Value = "'3.628"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "'121.41"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("D48").Value = "'1.893"
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").Value = "'0.06724"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").Value = "'70.84"
$ws.Range("E51").Value = "  +2.48%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.234"
$ws.Range("E40").Value = "  +10.19%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6121"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("D42").Value = "'1.428"
$ws.Range("D43").Value = "'1.134"
$ws.Range("E43").Value = "  +4.66%  "
